$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.352.12"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.593.30"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.816.49"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.607.39"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.356.52"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "1.339.49"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  -18.28%  "
$ws.Range("E41").Value = "  +5.25%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0985"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("E51").Value = "  -0.59%  "
